$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column M header (job_id), reusing the existing header style (s=1)
# by copying A1's format onto M1.
$ws.Range("M1").Value = "job_id"
$ws.Range("A1").Copy()
$ws.Range("M1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Row 46 ----
$ws.Range("A46").Value = 'schonherz'
$ws.Range("B46").Value = 'Web Developer Intern - Calibration Team'
$ws.Range("C46").Value = 'https://schonherz.hu/diakmunka/budapest/fejleszto---tesztelo/40918-web-developer-intern-calibration-team'
$ws.Range("D46").Value = '🔐 Work on real autonomous driving projects as a Web Developer Intern in a friendly, small team. You’ll build data-driven web applications for sensor calibration and analytics while gaining hands-on experience in automotive R&D with flexible hours and long-term opportunities. 🚗 Our Partner is a software development company specializing in autonomous vehicles. Their primary approach to self-driving technology is camera-based, relying on image processing, while also experimenting with radar and lidar solutions. The company is Hungarian-owned, yet operates with a distinctly American-style, relaxed management culture, works on numerous U.S.-based projects, and offers internal career advancement opportunities. 🪂 Developing automated driving requires a complete and mature toolchain to collect, generate, use, and manage the data needed for a safe and robust solution. A modern, data-driven pipeline for automated driving reduces the complexity of data processing through a high level of automation, while still ensuring automotive-grade quality. Sensors are the most important link between the system and its environment. Precise sensor calibration and a deep understanding of how sensors perceive the world are essential for self-driving applications. As a web developer intern in a small team, the intern takes part in the development of calibration and diagnostic data analytics pipelines and supports the team’s daily work through coding and testing. Job Summary 🌀 Collaborating closely with research colleagues to create visualizations for cutting-edge sensor calibration algorithms (React.js, Python) 🌀 Contribute to the development and testing of the calibration and diagnostic pipeline by building intuitive web applications to enhance data comprehension (JavaScript, Python) 🌀 Collaborate with users and the design team to enhance the user experience 🌀 Ensure the functionality and upkeep of the calibration and diagnostic database through thorough testing and maintenance (SQL) 🌀 Assist the team in developing and maintaining an automated testing pipeline, as well as supporting manual testing efforts Expectations: 🟨 Ongoing (active/passive) full-time university studies for at least one more year 🟨 Solid Python and React.js knowledge 🟨 Intermediate level English knowledge 🟨 Availability of at least 20 working hours per week 🟨 Ability to work independently What our Partner offers: 📍 Possibility to gain valuable work experience at the forefront of the automotive R&D sector 💸 Competitive student salary 💡 Inspiring and supportive working environment 🎳 Teambuilding events and other benefits for students ⏳ Flexible working hours 📗 Long-term job opportunity 🫟 Possibility of an immediate start 🧷 In-service training opportunity 🏋️ Fitness opportunities 🥯 Free lunch every workday at one of the best all-you-can-eat restaurants in Budapest or delivered by Wolt for Work ⏰ Working hours: From 8:00 a.m. to 6:00 p.m. ❗ How to apply: 1️⃣ Register on Schönherz.hu 2️⃣ Fill in a CV on Schönherz.hu * If you are over 25 years of age, 15% of your gross salary will be deducted from your PIT. If you are under 25 years old, you will be exempt from PIT, so you can consider your gross salary as net, up to the average wage in the national economy, which is HUF 693,700 per month.'
$ws.Range("E46").Value = 'English (Intermediate)'
$ws.Range("F46").Value = 'Budapest, District II. + Home Office'
$ws.Range("G46").Value = 'Gross 2250-2500 Ft/hour*'
$ws.Range("I46").Value = 'ml_score=3, ml_prob=0.28012610847425057'
$ws.Range("J46").Value = 3
$ws.Range("K46").Value = 0.2801261084742506
$ws.Range("L46").Value = 3
$ws.Range("M46").Value = 1

# ---- Row 47 ----
$ws.Range("A47").Value = 'muisz'
$ws.Range("B47").Value = 'Test Engineer (Intern) on Industrial Automation Field'
$ws.Range("C47").Value = 'https://muisz.hu/diakmunkaink/test-engineer-%20Industrial-Automation-Field-trainee-xi-kerulet--28375'
$ws.Range("D47").Value = 'Test Engineer (Intern) on Industrial Automation Field Today the community of around 1800 IT professionals and engineers believe - along with me - that the evosoft team is the place for us. Our co-workers in Budapest, Miskolc and Szeged are busy writing software which, although invisible to the naked eye, can fundamentally influence the very basics of our everyday lives. You will come across our software solutions in different kinds of the largest medical equipment technology, not to mention that our codes drive the greatest car manufacturers’ automation systems and you can even find them/us in electric cars. There are hardly any industrial areas where the programs I and my colleagues tested or developed are not present. The company evosoft, in its full name evosoft Hungary Kft. (Ltd.) started with 3 members more than 30 years ago. As of today, our headcount has reached the 1800 employees and we are expanding - join us, be part of our team!" Would you like to get to know those Siemens devices, which control the world’s largest and most innovative factories? Would you like to be part of a professional team, where you can constantly challenge yourself in your daily work? Would you like to work in an attractive lab environment, where you can constantly develop yourself as part of a dynamic team? If your answer is YES, then read our offer below! Feladatok •	Working closely with Test Engineers and Developers to understand product requirements and ensure high quality in releases •	Testing Siemens SIMATIC S7 PLC and HMI product families’ engineering systems and hardwares •	Executing test cases (manual and/or automated) and analyzing results •	Writing and maintaining test plans / test cases •	Designing and building hardware configurations for tests •	Bug analysis, documentation and follow up in the bug documenting system •	Participation in the setup of the test environment Elvárások •	Ongoing higher education (electrical engineer, software engineer, mechanical engineer, mechatronics engineer or similar) •	Basic electrical knowledge •	Knowledge of MS Windows 11, MS Office •	Communication level English language skills •	Ability to work independently and in a team, flexibility •	Independent, able to solve problems, analytical mindset •	Motivation to continuous self-development Egyéb Advantages: •	PLC programming knowledge (e.g. STL / SCL / FBD / LAD) •	Knowledge of Siemens SIMATIC S7-300 / 400 / 1200 / 1500 PLC product families •	Knowledge of Siemens STEP7 / TIA Portal product families •	C# and/or Python programming knowledge What we offer: •	Real-world experience in agile development processes •	You can work on various challenging tasks in the field of industrial automation •	Flexible working hours and the option to work remotely •	A dynamic, supporting team with mentor(s) •	Attractive working environment and secure workplace •	Free coffee and fruit days •	Recreational opportunities within the office building (table tennis, darts etc.) •	Various career opportunities •	Joining possibility to our different communities: sports communities, CSR, green, ToastMasters etc. •	Parking lot within the office building'
$ws.Range("E47").Value = 'magyar, angol'
$ws.Range("F47").Value = 'Budapest 11. ker.'
$ws.Range("G47").Value = '2200-2800 Ft/óra'
$ws.Range("I47").Value = 'ml_score=3, ml_prob=0.2318041591459815'
$ws.Range("J47").Value = 3
$ws.Range("K47").Value = 0.2318041591459815
$ws.Range("L47").Value = 3
$ws.Range("M47").Value = 2

# ---- Row 48 ----
$ws.Range("A48").Value = 'muisz'
$ws.Range("B48").Value = 'Android developer trainee'
$ws.Range("C48").Value = 'https://muisz.hu/diakmunkaink/Android-developer-trainee-III-kerulet--27760'
$ws.Range("D48").Value = 'Android developer trainee Android developer trainee Feladatok The main task of the candidate would be to write unit tests for already completed components, for which it would be useful to use AI assistance. What would be your responsibilities? •	Develop HMI components for automotive cockpit systems •	Implement unit tests and other improvements for components •	Support integration with vehicle middleware and hardware on Android Automotive OS Elvárások The ideal candidate should possess experience in the following •	Kotlin development •	Android HMI •	MVVM pattern and unit testing •	AI assistance in software development •	MS Teams, Jira, Confluence and Bitbucket/Git •	Good verbal and written English Egyéb - Participate in team meetings and brainstorming sessions, contributing ideas for new marketing campaigns and content strategies. - Assist with the creation of marketing materials, including brochures, flyers, and presentations, as needed. - Perform administrative tasks related to marketing projects and campaigns, ensuring smooth operations within the team.'
$ws.Range("E48").Value = 'magyar, angol'
$ws.Range("F48").Value = 'Budapest 3. ker.'
$ws.Range("G48").Value = '2600 Ft/óra'
$ws.Range("I48").Value = 'ml_score=3, ml_prob=0.2557808943824647'
$ws.Range("J48").Value = 3
$ws.Range("K48").Value = 0.2557808943824647
$ws.Range("L48").Value = 3
$ws.Range("M48").Value = 3
